$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.72093023255813948
$ws.Range("C2").Value = [double]"6.1353487706492855E-2"
$ws.Range("D2").Value = 2
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = [double]"1.8632452769907909E-2"
$ws.Range("D3").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = [double]"1.4661379753313216E-2"
$ws.Range("D4").Value = 20
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = [double]"9.0227223098169121E-3"
$ws.Range("D5").Value = 3
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = [double]"1.4661123753247681E-2"
$ws.Range("D6").Value = 7
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = [double]"7.744769982661116E-3"
$ws.Range("D7").Value = 5
$ws.Range("B8").Value = 0.72093023255813948
$ws.Range("C8").Value = [double]"7.5051539213194041E-3"
$ws.Range("D8").Value = 19
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = [double]"7.2931858670555816E-3"
$ws.Range("D9").Value = 11
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = [double]"7.6764179651629995E-3"
$ws.Range("D10").Value = 16
$ws.Range("B11").Value = 0.7441860465116279
$ws.Range("C11").Value = [double]"9.9184665391274345E-3"
$ws.Range("D11").Value = 13
$ws.Range("C12").Value = [double]"4.4149003302144843E-2"
$ws.Range("D12").Value = 13
$ws.Range("C13").Value = [double]"8.9960983030011651E-3"
$ws.Range("D13").Value = 12
$ws.Range("C14").Value = [double]"8.9909783016904454E-3"
$ws.Range("D14").Value = 9
$ws.Range("C15").Value = [double]"7.4864659165352751E-3"
$ws.Range("D15").Value = 10
$ws.Range("C16").Value = [double]"7.2883218658103977E-3"
$ws.Range("D16").Value = 5
$ws.Range("B17").Value = 0.70930232558139539
$ws.Range("C17").Value = [double]"7.7396499813503954E-3"
$ws.Range("D17").Value = 15
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = [double]"7.273985862140381E-3"
$ws.Range("D18").Value = 3
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = [double]"7.0891538148233763E-3"
$ws.Range("D19").Value = 20
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = [double]"7.7667859882972132E-3"
$ws.Range("D20").Value = 14
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = [double]"1.0642178724397754E-2"
$ws.Range("D21").Value = 18
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = [double]"7.9208980277498953E-3"
$ws.Range("B23").Value = 0.73255813953488369
$ws.Range("C23").Value = [double]"7.3489938813424338E-3"
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = [double]"7.3958418933355247E-3"
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = [double]"7.4572819090641691E-3"
$ws.Range("B26").Value = 0.7558139534883721
$ws.Range("C26").Value = [double]"7.5328019283972935E-3"
$ws.Range("C27").Value = [double]"7.5458579317396303E-3"
$ws.Range("C28").Value = [double]"7.237889852899802E-3"
$ws.Range("C29").Value = [double]"7.4186258991682302E-3"
$ws.Range("C30").Value = [double]"1.1177218861368028E-2"
$ws.Range("C31").Value = [double]"7.3000978688250546E-3"
$ws.Range("B32").Value = 0.82558139534883723
$ws.Range("C32").Value = [double]"1.1446274930246383E-2"
$ws.Range("B33").Value = 0.84705882352941175
$ws.Range("C33").Value = [double]"7.53689792944587E-3"
$ws.Range("B34").Value = 0.76470588235294112
$ws.Range("C34").Value = [double]"7.2862738652861095E-3"
$ws.Range("C35").Value = [double]"7.2793618635166374E-3"
$ws.Range("B36").Value = 0.82352941176470584
$ws.Range("C36").Value = [double]"7.3167378730848954E-3"
$ws.Range("B37").Value = 0.81395348837209303
$ws.Range("C37").Value = [double]"7.3116178717741748E-3"
$ws.Range("B38").Value = 0.88372093023255816
$ws.Range("C38").Value = [double]"7.1518738308797005E-3"
$ws.Range("B39").Value = 0.86046511627906974
$ws.Range("C39").Value = [double]"7.1892498404479595E-3"
$ws.Range("B40").Value = 0.81395348837209303
$ws.Range("C40").Value = [double]"7.2806418638443173E-3"
$ws.Range("B41").Value = 0.84883720930232553
$ws.Range("C41").Value = [double]"7.2450578547348108E-3"
$ws.Range("B47").Value = 0.81395348837209303
$ws.Range("B48").Value = 0.88235294117647056
$ws.Range("B49").Value = 0.82352941176470584
$ws.Range("B50").Value = 0.90588235294117647
$ws.Range("B51").Value = 0.83529411764705885
$ws.Range("B52").Value = 0.80232558139534882
$ws.Range("B54").Value = 0.89534883720930236
$ws.Range("B55").Value = 0.83720930232558144
$ws.Range("B56").Value = 0.87209302325581395
